# Applies the coin-price refresh captured in the commit diff.
# D-column cells whose new text reads as a plain decimal (e.g. "89.55") are
# pre-formatted as Text ("@") so Excel stores the literal digit string instead
# of silently parsing it into a float (which would also eat trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.643.18'
$ws.Range("E2").Value = '  +1.31%  '
# Row 3
$ws.Range("D3").Value = '2.233.87'
$ws.Range("E3").Value = '  -0.12%  '
# Row 4
$ws.Range("E4").Value = '  -0.05%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.65'
$ws.Range("E5").Value = '  +2.82%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '89.55'
# Row 7
$ws.Range("E7").Value = '  +0.89%  '
# Row 8
$ws.Range("E8").Value = '  -0.03%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  +1.03%  '
# Row 10
$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.61'
$ws.Range("E10").Value = '  +4.13%  '
# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.80'
$ws.Range("E11").Value = '  +7.48%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0788'
$ws.Range("E12").Value = '  +0.94%  '
# Row 13
$ws.Range("E13").Value = '  +2.69%  '
# Row 14
$ws.Range("E14").Value = '  -0.01%  '
# Row 15
$ws.Range("D15").Value = '2.578.41'
$ws.Range("E15").Value = '  +0.03%  '
# Row 16
$ws.Range("E16").Value = '  +1.02%  '
# Row 17
$ws.Range("D17").Value = '2.241.71'
$ws.Range("E17").Value = '  +0.64%  '
# Row 18
$ws.Range("E18").Value = '  +1.93%  '
# Row 19
$ws.Range("D19").Value = '40.507.23'
$ws.Range("E19").Value = '  +1.20%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0896'
$ws.Range("E20").Value = '  +0.56%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.51'
$ws.Range("E21").Value = '  +0.81%  '
# Row 22
$ws.Range("E22").Value = '  -0.05%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.04'
$ws.Range("E23").Value = '  +0.63%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.37'
$ws.Range("E24").Value = '  +0.24%  '
# Row 25
$ws.Range("E25").Value = '  +2.70%  '
# Row 26
$ws.Range("E26").Value = '  -0.14%  '
# Row 27
$ws.Range("E27").Value = '  +1.93%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.82'
$ws.Range("E28").Value = '  +4.50%  '
# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.44'
$ws.Range("E29").Value = '  +2.27%  '
# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.09'
$ws.Range("E30").Value = '  -4.63%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.04'
$ws.Range("E31").Value = '  +0.77%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.89'
$ws.Range("E32").Value = '  +2.85%  '
# Row 33
$ws.Range("E33").Value = '  +0.10%  '
# Row 34
$ws.Range("E34").Value = '  +2.18%  '
# Row 35
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.03'
$ws.Range("E35").Value = '  +5.07%  '
# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0725'
$ws.Range("E36").Value = '  +0.95%  '
# Row 37
$ws.Range("E37").Value = '  -0.34%  '
# Row 38
$ws.Range("E38").Value = '  +6.55%  '
# Row 39
$ws.Range("E39").Value = '  +1.86%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("E40").Value = '  +4.48%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.95'
$ws.Range("E41").Value = '  +1.00%  '
# Row 42
$ws.Range("E42").Value = '  -0.15%  '
# Row 43
$ws.Range("D43").Value = '2.091.44'
$ws.Range("E43").Value = '  -1.67%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.64'
$ws.Range("E44").Value = '  +7.08%  '
# Row 45
$ws.Range("E45").Value = '  +2.36%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.12'
$ws.Range("E46").Value = '  +2.63%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.87'
$ws.Range("E47").Value = '  +7.68%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.85'
$ws.Range("E48").Value = '  -13.06%  '
# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.51'
$ws.Range("E49").Value = '  +2.24%  '
# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.447.92'
$ws.Range("E50").Value = '  +0.29%  '
# Row 51
$ws.Range("E51").Value = '  +2.55%  '
